# "Adding comment Test1 in OR"
# - Remove the "system" worksheet entirely.
# - Repurpose the "testdata" worksheet into a new "screenTitles" worksheet
#   that maps screen objectIDs to their Dutch display titles.
# - Update the "hub" sheet's view (selection moves to A7, it is no longer
#   the active/selected tab - "screenTitles" becomes the active tab instead).

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$hub = $wb.Worksheets.Item("hub")
$testdata = $wb.Worksheets.Item("testdata")
$system = $wb.Worksheets.Item("system")

# Remove the obsolete "system" sheet - its content moved into "screenTitles".
$system.Delete()

# Turn the old "testdata" (username/password) sheet into the new
# "screenTitles" lookup sheet.
$testdata.Cells.ClearContents()
$testdata.Name = "screenTitles"

# Header row - copy the bold/shaded header style used on "hub" so the new
# header cells share the same cell style (s="2") instead of minting a new one.
$hub.Range("A1").Copy()
$testdata.Range("A1:B1").PasteSpecial(-4122)

$testdata.Range("A1").Value = "objectID"
$testdata.Range("B1").Value = "name_nl"

$testdata.Range("B4").Value = "epg"
$testdata.Range("B2").Value = "instellingen"
$testdata.Range("A3").Value = "System"
$testdata.Range("A2").Value = "Setting"
$testdata.Range("B3").Value = "systeem"
$testdata.Range("A4").Value = "epgSetting"

$testdata.Columns.Item(1).ColumnWidth = 14.67
$testdata.Columns.Item(2).ColumnWidth = 17.83

# "hub" keeps its data untouched; only the view/selection changes - it is no
# longer the selected tab.
$hub.Range("A7").Select()

# "screenTitles" becomes the active tab, selection moves to T25.
$testdata.Activate()
$testdata.Range("T25").Select()
